$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 369 (shifts 369..388 down to 370..389)
$ws.Rows.Item(369).Insert()

# Populate the new row 369 with the new record
$ws.Range("A369").Value = 4
$ws.Range("B369").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C369").Value = 'Los Lagos'
$ws.Range("D369").Value = 44931
$ws.Range("E369").Value = 10
$ws.Range("F369").Value = 100112003
$ws.Range("G369").Value = 'Ajo'
$ws.Range("H369").Value = 'Chino'
$ws.Range("I369").Value = 'Primera'
$ws.Range("J369").Value = 80
$ws.Range("K369").Value = 18000
$ws.Range("L369").Value = 18000
$ws.Range("M369").Value = 18000
$ws.Range("N369").Value = '$/caja 10 kilos'
$ws.Range("O369").Value = 'China'
$ws.Range("P369").Value = 1800
$ws.Range("Q369").Value = 10
$ws.Range("R369").Value = 'Hortaliza'
